$d = $word.ActiveDocument

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$wdHeaderFooterPrimary = 1
$wdHeaderFooterFirstPage = 2

foreach ($sec in $d.Sections) {

    # Footer logos (Pearson logo): image2.png -> image1.png
    foreach ($idx in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage)) {
        $ftr = $sec.Footers.Item($idx)
        if ($ftr.Exists) {
            for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
                $shp = $ftr.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }

    # Header logos (BTEC logo): image1.jpg -> image2.jpg
    foreach ($idx in @($wdHeaderFooterPrimary, $wdHeaderFooterFirstPage)) {
        $hdr = $sec.Headers.Item($idx)
        if ($hdr.Exists) {
            for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
                $shp = $hdr.Range.InlineShapes.Item($i)
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
}
